$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @"
2025-08-21T07:22:57|Haus_B.3dm|0525f66e-fe44-413e-a4d2-beb390edd8e8|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:23:34|Haus_B.3dm|12bb73f7-d6f1-4f90-b9c4-d824137b5ca1|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:24:10|Haus_B.3dm|7de5608d-38d2-4592-9d92-4ee321e99e57|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:25:40|Haus_B.3dm|436612dc-8e28-42e1-b7b9-9d8ddd7d1723|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:25:55|Haus_B.3dm|8c65d0d3-12f6-47bb-ad27-374236345c7f|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:26:36|Haus_B.3dm|a24634b9-c311-4fd4-a90b-2056fbd59603|rahmentuere_w|Standard 1:10 Rahmenbeschriftung WHG Eingang|rahmentuerew.csv
2025-08-21T07:27:17|Haus_B.3dm|a4942db9-b96b-4091-9a84-fc7a398957ad|rahmentuere_w|Standard 1:10 Rahmenbeschriftung WHG Eingang|rahmentuerew.csv
2025-08-21T07:28:05|Haus_B.3dm|c3a4f394-8a16-4450-a9ce-e7ea79494d13|rahmentuere_w|Standard 1:10 Rahmenbeschriftung WHG Eingang|rahmentuerew.csv
2025-08-21T07:28:30|Haus_B.3dm|9dd27323-a1cb-4594-bfd8-82b762ca7d4f|rahmentuere|Standard 1:10 Rahmenbeschriftung|rahmentuere.csv
2025-08-21T07:34:57|Haus_B.3dm|5ec54a75-e768-424a-929d-08a6a2f585d9|rahmentuere|Standard 1:10 Rahmenbeschriftung|rahmentuere.csv
2025-08-21T07:35:55|Haus_B.3dm|0c0c8aa3-cdc9-4cbc-a3c8-16510966b3b9|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:36:42|Haus_B.3dm|31ef67a2-3823-4328-a160-fa8fb78b95b9|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:37:26|Haus_B.3dm|0ef27e03-7055-4f17-a320-3787d17657b8|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:37:48|Haus_B.3dm|f22e8e25-09db-4697-936f-ed849427c2f6|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:38:19|Haus_B.3dm|78f4d977-3f50-4d89-9526-85cbbdefb145|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:38:45|Haus_B.3dm|b870949e-ab73-4dc5-8cd6-a190aa69b45f|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:39:09|Haus_B.3dm|e9743612-ad46-41d9-807d-843f5efff1e1|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:39:31|Haus_B.3dm|afeb0386-c654-46b3-b397-b978512f1737|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T07:41:08|Haus_B.3dm|64d94098-810c-4b12-88bd-10c1e198a48d|schiebetuere|Standard 1:10 Schiebetürbeschriftung|schiebetuere.csv
2025-08-21T07:41:37|Haus_B.3dm|7e039f2e-35fe-4dbd-8e09-06ff49750665|schiebetuere|Standard 1:10 Schiebetürbeschriftung|schiebetuere.csv
2025-08-21T15:20:04|Haus_A.3dm|1c202ef5-3c9c-4e35-a924-d2ae95d019ab|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T15:33:31|Haus_A.3dm|a9760fd2-627e-45ab-bbcc-54109be14cc5|schiebetuere|Standard 1:10 Schiebetürbeschriftung|schiebetuere.csv
2025-08-21T15:34:37|Haus_A.3dm|ef395ba8-6040-4c51-8d1e-eb597ff5158b|spez|Standard 1:10 Spez.Rahmenbeschriftung|spez.csv
2025-08-21T15:40:32|Haus_A.3dm|92869a40-4619-440f-b718-352db991f927|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T16:10:20|Haus_D_test.3dm|7258b958-748a-4a7a-b9dc-fb1dc1284e61|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T16:29:06|Haus_D_test.3dm|12305b9a-7c5e-4a5c-81bc-9016e3fdd0d0|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere_2.csv
2025-08-21T16:34:27|Haus_D_test.3dm|6f11d915-9e6a-4c02-90a2-a8ff763f6018|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere_2.csv
2025-08-21T16:35:05|Haus_D_test.3dm|abf55ca4-1ec3-4367-b8af-491d9cc13913|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere.csv
2025-08-21T16:35:22|Haus_D_test.3dm|d8cc6915-2784-4c94-bb37-c6945dc6e50f|rahmentuere|Standard 1:10 Rahmenbeschriftung|rahmentuere.csv
2025-08-21T16:39:10|Haus_D_test.3dm|c93fd141-acf6-49a1-9aa0-9cffd06b00ec|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere_2.csv
2025-08-21T16:40:39|Haus_D_test.3dm|86a7fca8-f4f1-4248-ac47-aa6359fc762c|rahmentuere|Standard 1:10 Rahmenbeschriftung|rahmentuere.csv
2025-08-21T16:52:59|Haus_D_test.3dm|96587c5c-6bcf-4357-80ad-b6b43e7414b1|zargentuere|Standard 1:10 Zargenbeschriftung|zargentuere_2.csv
"@

$lines = $rowData -split "`n"
$startRow = 30

$r = $startRow
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $fields = $line -split '\|'
    $ws.Cells.Item($r, 1).Value = $fields[0]
    $ws.Cells.Item($r, 2).Value = $fields[1]
    $ws.Cells.Item($r, 3).Value = $fields[2]
    $ws.Cells.Item($r, 4).Value = $fields[3]
    $ws.Cells.Item($r, 5).Value = $fields[4]
    $ws.Cells.Item($r, 6).Value = $fields[5]
    $r++
}
